$wb = $excel.ActiveWorkbook

# --- Update "Logs" sheet: append new row 9 ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A9").Value = "Bel jij klant Jansen even?"
$logs.Range("B9").Value = "mailmind.test@zohomail.eu"
$logs.Range("C9").Value = "Testmail #19: Bel jij klant Jansen even?"
$logs.Range("D9").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("E9").Value = "Beste,`nBedankt voor je e-mail. Het spijt me, maar ik heb geen verdere informatie over welke klant Jansen je bedoelt. Zou je meer context kunnen geven, zodat ik je beter kan helpen?`nMet vriendelijke groet,`n[Je naam]`nE-mailassistent"
$logs.Range("F9").Value = "2025-08-02 00:16:41"
$logs.Range("G9").Value = "Ja"
$logs.Range("H9").Value = "Nee"
$logs.Range("I9").Value = "Ja"
$logs.Range("J9").Value = "Nee"

# Undo the automatic row-height change caused by the long multi-line E9 text
# so the row matches the other (non custom-height) rows.
$logs.Rows.Item(9).AutoFit()

# --- Extend conditional formatting ranges to include the new row ---
foreach ($col in @("D","G","H","I","J")) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "8")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "9")
    $fcs = $oldRange.FormatConditions
    $count = $fcs.Count
    for ($i = 1; $i -le $count; $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($newRange)
    }
}

# --- Update "Dashboard" sheet: swap category rows 2 and 3 (counts of the ---
# --- "Intern verzoek" category increases to 3 after the new row above)  ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A2").Value = "Intern verzoek / Actie voor medewerker"
$dash.Range("B2").Value = 3
$dash.Range("A3").Value = "Productinformatie"
$dash.Range("B3").Value = 2
